$d = $word.ActiveDocument

$targets = @(
    "Songs must be playable on both Windows and MACos",
    "Next and previous marker commands should work also when the song in paused.",
    "Song position should return to zero when the song reaches the end."
)

foreach ($paragraph in $d.Paragraphs) {
    $text = $paragraph.Range.Text.TrimEnd([char]13, [char]7)
    foreach ($target in $targets) {
        if ($text -eq $target) {
            $paragraph.Range.Font.StrikeThrough = 1
        }
    }
}
